# Rename the two BTec logo pictures (one in each header variant) and the
# two Pearson logo pictures (one in each footer variant) - swapping their
# cosmetic "image1.jpg"/"image2.jpg" and "image1.png"/"image2.png" names,
# exactly as described by the commit diff. Only the picture's display
# Name changes; size, position, alt-text and the underlying image data
# are left untouched.

$d = $word.ActiveDocument
$sec = $d.Sections.First

# --- Headers: BTec_Logo-Orange, "image2.jpg" -> "image1.jpg" ---------------

# Primary (default) header
$headerDefault = $sec.Headers.Item(1)
$btecDefault = $headerDefault.Range.InlineShapes.Item(1)
if ($btecDefault.AlternativeText -eq "BTec_Logo-Orange") {
    $btecDefault.Name = "image1.jpg"
}

# First-page header
$headerFirst = $sec.Headers.Item(2)
$btecFirst = $headerFirst.Range.InlineShapes.Item(1)
if ($btecFirst.AlternativeText -eq "BTec_Logo-Orange") {
    $btecFirst.Name = "image1.jpg"
}

# --- Footers: PearsonLogo.png, "image1.png" -> "image2.png" ----------------

# Primary (default) footer
$footerDefault = $sec.Footers.Item(1)
$pearsonDefault = $footerDefault.Range.InlineShapes.Item(1)
if ($pearsonDefault.AlternativeText -like "*PearsonLogo.png") {
    $pearsonDefault.Name = "image2.png"
}

# First-page footer
$footerFirst = $sec.Footers.Item(2)
$pearsonFirst = $footerFirst.Range.InlineShapes.Item(1)
if ($pearsonFirst.AlternativeText -like "*PearsonLogo.png") {
    $pearsonFirst.Name = "image2.png"
}
